$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record (week) was inserted into the daily price log.
# This shifts all existing rows from 63 downward by one row (the
# previous row 135 now becomes row 136), and the newly opened row 63
# is populated with the new week's data.
$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value2 = 9
$ws.Range("B63").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value2 = "Metropolitana"
$ws.Range("D63").Value2 = 45175
$ws.Range("E63").Value2 = 13
$ws.Range("F63").Value2 = 100112005
$ws.Range("G63").Value2 = "Puerro"
$ws.Range("H63").Value2 = "Sin especificar"
$ws.Range("I63").Value2 = "Primera"
$ws.Range("J63").Value2 = 70
$ws.Range("K63").Value2 = 8000
$ws.Range("L63").Value2 = 8000
$ws.Range("M63").Value2 = 8000
$ws.Range("N63").Value2 = "`$/paquete 20 unidades"
$ws.Range("O63").Value2 = "Provincia de Chacabuco"
$ws.Range("P63").Value2 = 400
$ws.Range("Q63").Value2 = 20
$ws.Range("R63").Value2 = "Hortaliza"
